$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data from the latest GitHub Actions run.
# Price (column D) and Volume(1h) (column E) values are stored as plain text in the
# source sheet, so purely-numeric-looking Price values are prefixed with a leading
# apostrophe to force Excel to keep them as text instead of auto-converting to numbers.
$ws.Range("D2").Value = "55.017.62"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "2.296.30"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'506.86"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").Value = "'130.05"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.532"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").Value = "2.319.02"
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("E10").Value = "  +2.66%  "
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("D12").Value = "'5.08"
$ws.Range("E12").Value = "  +6.95%  "
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "'23.89"
$ws.Range("E14").Value = "  +4.41%  "
$ws.Range("D15").Value = "2.705.78"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "55.039.87"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").Value = "2.369.27"
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("D19").Value = "'10.60"
$ws.Range("E19").Value = "  +3.31%  "
$ws.Range("D20").Value = "'4.19"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").Value = "'310.96"
$ws.Range("E21").Value = "  +2.39%  "
$ws.Range("D22").Value = "'6.62"
$ws.Range("E22").Value = "  +4.01%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "'60.55"
$ws.Range("E24").Value = "  -2.19%  "
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'7.52"
$ws.Range("E27").Value = "  +2.36%  "
$ws.Range("D28").Value = "'173.13"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("E29").Value = "  +3.05%  "
$ws.Range("D30").Value = "0.0₃0709"
$ws.Range("E30").Value = "  +2.19%  "
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("E32").Value = "  +5.45%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'18.06"
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("E36").Value = "  +2.74%  "
$ws.Range("D37").Value = "'0.919"
$ws.Range("E37").Value = "  -4.87%  "
$ws.Range("E38").Value = "  +4.24%  "
$ws.Range("D39").Value = "'36.84"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("E40").Value = "  +2.29%  "
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("D42").Value = "'134.71"
$ws.Range("E42").Value = "  +6.68%  "
$ws.Range("D43").Value = "'3.44"
$ws.Range("E44").Value = "  -2.28%  "
$ws.Range("D45").Value = "'261.71"
$ws.Range("E45").Value = "  +7.95%  "
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("D49").Value = "'0.377"
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("D51").Value = "'16.51"
$ws.Range("E51").Value = "  +0.70%  "
